# Update the published output of the StructureDefinition-mindfulness-audit-level
# spreadsheet for the "2025 august" refresh:
#   - Metadata!URL / Elements!Extension.url fixed value -> new canonical base URL
#   - Metadata!Date -> new generation timestamp
#   - Elements! Binding Value Set -> new ValueSet URL

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

$newUrl = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/mindfulness-audit-level"
$newDate = "2025-08-20T10:40:04+01:00"
$newValueSet = "https://2rdoc.pt/fhir/ValueSet/audit-levels"

# Metadata sheet: URL (row 2) and Date (row 8) values live in column B.
$wsMetadata.Range("B2").Value = $newUrl
$wsMetadata.Range("B8").Value = $newDate

# Elements sheet: the Extension.url row's "Fixed Value" column (R) repeats the
# canonical URL, and the last slice's "Binding Value Set" column (Z) points at
# the audit-levels value set.
$wsElements.Range("R5").Value = $newUrl
$wsElements.Range("Z6").Value = $newValueSet
